# Weekly update: insert a new price record as row 48 (Feria Lagunitas de
# Puerto Montt - Espárragos), pushing the existing rows 48-69 down to 49-70.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 48, shifting rows 48:69 down
# to 49:70 (dimension grows from A1:R69 to A1:R70 automatically).
$ws.Rows.Item(48).Insert()

# Fill in the new record in row 48.
$ws.Range("A48").Value2 = 4
$ws.Range("B48").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C48").Value2 = "Los Lagos"
$ws.Range("D48").Value2 = 45202
$ws.Range("E48").Value2 = 10
$ws.Range("F48").Value2 = 300000000
$ws.Range("G48").Value2 = "Espárragos"
$ws.Range("H48").Value2 = "Sin especificar"
$ws.Range("I48").Value2 = "Primera"
$ws.Range("J48").Value2 = 600
$ws.Range("K48").Value2 = 1900
$ws.Range("L48").Value2 = 2300
$ws.Range("M48").Value2 = 2100
$ws.Range("N48").Value2 = "`$/kilo"
$ws.Range("O48").Value2 = "Provincia de Linares"
$ws.Range("P48").Value2 = 2100
$ws.Range("Q48").Value2 = 1
$ws.Range("R48").Value2 = "Hortaliza"

# Match the date-formatted style used by the other rows in column D.
$ws.Range("D48").NumberFormat = $ws.Range("D49").NumberFormat
